# "add stats vs BGYfk"
# - Refresh the Д3 standings table (rows 5-18) with the latest games/wins/
#   losses/points/score totals, which also re-sorts several teams in the
#   table.
# - Append the newest two match days (15 and 16 Feb 2025 -> serials 45703
#   and 45704) with their results to the results log (rows 102-110).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Д3")

# ---------------------------------------------------------------------
# 1. Standings table (B4:H18) - update team order, record and points.
#    Columns: B=place, C=team, D=games, E=wins, F=losses, G=score, H=points
# ---------------------------------------------------------------------
$standings = @(
    @(1,  "ISsoft",                 10, 8, 2, "653 - 532", 18),
    @(2,  "Эра-Недвижимости плюс",  10, 8, 2, "763 - 601", 18),
    @(3,  "ОПЛАТИ",                 10, 8, 2, "756 - 616", 18),
    @(4,  "GOLDEN HILL",            10, 8, 2, "728 - 668", 18),
    @(5,  "Грушвиль",               10, 7, 3, "805 - 657", 17),
    @(6,  "БГУФК",                  10, 7, 3, "694 - 545", 17),
    @(7,  "Mapogo males",           10, 7, 3, "750 - 694", 17),
    @(8,  "SIRIUS",                 10, 5, 5, "677 - 588", 15),
    @(9,  "Стрела",                 10, 3, 7, "600 - 655", 13),
    @(10, "Eagles",                 10, 3, 7, "578 - 619", 13),
    @(11, "VSS",                    10, 3, 7, "621 - 688", 13),
    @(12, "NORD",                   10, 2, 8, "508 - 791", 12),
    @(13, "ЛФК",                    10, 1, 9, "551 - 730", 11),
    @(14, "Минск 7х",               10, 0, 10, "437 - 737", 10)
)

$row = 5
foreach ($team in $standings) {
    $ws.Cells.Item($row, 2).Value = $team[0]   # B - place
    $ws.Cells.Item($row, 3).Value = $team[1]   # C - team name
    $ws.Cells.Item($row, 4).Value = $team[2]   # D - games
    $ws.Cells.Item($row, 5).Value = $team[3]   # E - wins
    $ws.Cells.Item($row, 6).Value = $team[4]   # F - losses
    $ws.Cells.Item($row, 7).Value = $team[5]   # G - score
    $ws.Cells.Item($row, 8).Value = $team[6]   # H - points
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Append new match-day blocks to the results log.
#    Row pattern: a shaded date row followed by one merged row per game.
# ---------------------------------------------------------------------
$lastRow = 101

$day1Date = 45703
$day1Games = @(
    "ISsoft - Eagles 69:42 (16:30, БНТУ)",
    "ЛФК - Стрела 58:68 (18:00, БНТУ)",
    "VSS - Mapogo males 75:87 (19:30, БНТУ)"
)

$day2Date = 45704
$day2Games = @(
    "БГУФК - SIRIUS 59:56 (11:00, БНТУ)",
    "Грушвиль - GOLDEN HILL 80:85 (12:30, БНТУ)",
    "Эра-Недвижимости плюс - Минск 7х 92:43 (14:00, БНТУ)",
    "NORD - ОПЛАТИ 45:85 (15:30, БНТУ)"
)

$r = $lastRow

# date header row (style copied from the previous date row).
# NOTE: merge the (still empty) row *before* pasting the formatting -
# merging first avoids Excel re-deriving per-cell border edges for an
# existing/populated merge, which would otherwise fragment the style
# used by every other date/result row in the sheet.
$r = $r + 1
$ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).Merge()
$ws.Range("B93:H93").Copy()
$ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).PasteSpecial(-4122)
$ws.Cells.Item($r, 2).Value = $day1Date

foreach ($game in $day1Games) {
    $r = $r + 1
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).Merge()
    $ws.Range("B94:H94").Copy()
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).PasteSpecial(-4122)
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).RowHeight = 19.95
    $ws.Cells.Item($r, 2).Value = $game
}

# second date header row
$r = $r + 1
$ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).Merge()
$ws.Range("B93:H93").Copy()
$ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).PasteSpecial(-4122)
$ws.Cells.Item($r, 2).Value = $day2Date

foreach ($game in $day2Games) {
    $r = $r + 1
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).Merge()
    $ws.Range("B94:H94").Copy()
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).PasteSpecial(-4122)
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 8)).RowHeight = 19.95
    $ws.Cells.Item($r, 2).Value = $game
}

$excel.CutCopyMode = $false
